# "Generate Report for handoff"
# Bump the "Latest Handoff Datetime" (column D) for the most recently
# handed-off file (row 5: 6e6a413c-...md) on each language report sheet.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-28 08:49:13"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-28 08:49:25"
